function Find-ParaIndex {
    param($d, [string]$text)
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Replace-ParasWithXml {
    param($d, [int]$startIdx, [int]$endIdx, [string]$innerXml)
    if ($startIdx -lt 1 -or $endIdx -lt 1 -or $endIdx -lt $startIdx) {
        throw "Replace-ParasWithXml: invalid paragraph index range ($startIdx .. $endIdx)"
    }
    $pStart = $d.Paragraphs($startIdx)
    $pEnd = $d.Paragraphs($endIdx)
    $full = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $target = $d.Range($full.Start, $full.End - 1)
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$target.InsertXML($xmlFrag)
}

$d = $word.ActiveDocument

$PBDR = '<w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr>'

# ---- Hunk 2 (lower in doc, do first so hunk 1 indices aren't affected): "And here are two more pictures..." section ----
$startIdx = Find-ParaIndex $d "And here are two more pictures you’ll see in the game. How are these pictures different?"
$endIdx = Find-ParaIndex $d "Wiggle’s face is different in the two pictures."
$inner = '<w:p><w:pPr>' + $PBDR + '</w:pPr>' +
    '<w:r><w:t xml:space="preserve">And here are two more pictures you’ll see in the game. </w:t></w:r>' +
    '<w:r><w:t>What’s the difference between two pictures?</w:t></w:r>' +
    '</w:p>'
Replace-ParasWithXml $d $startIdx $endIdx $inner

# ---- Hunk 1: "To figure out..." through "The Wuggle that is talking..." section ----
$startIdx = Find-ParaIndex $d "To figure out what Wiggle will do next, you’ll have to use what you see in the picture."
$endIdx = Find-ParaIndex $d "The Wuggle that is talking to Wiggle is different in the two pictures."
$inner = '<w:p><w:pPr>' + $PBDR + '</w:pPr>' +
    '<w:r><w:t xml:space="preserve">Before you </w:t></w:r>' +
    '<w:r><w:t>start</w:t></w:r>' +
    '<w:r><w:t>, we’re going to show you some of the pictures</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> you’ll see in the game</w:t></w:r>' +
    '<w:r><w:t>. Your job is to tell us how the pictures are different.</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr>' + $PBDR + '</w:pPr></w:p>' +
    '<w:p><w:pPr>' + $PBDR + '</w:pPr>' +
    '<w:r><w:t xml:space="preserve">Here are two of the pictures you’ll see in the game. </w:t></w:r>' +
    '<w:r><w:t>What’s the difference between these pictures?</w:t></w:r>' +
    '</w:p>'
Replace-ParasWithXml $d $startIdx $endIdx $inner

# ---- Hunk 3: add lastRenderedPageBreak to "Wiggle feel zavy." paragraph ----
$idx = Find-ParaIndex $d "Wiggle feel zavy."
$inner = '<w:p><w:pPr>' + $PBDR + '</w:pPr>' +
    '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Wiggle feel </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>zavy</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
Replace-ParasWithXml $d $idx $idx $inner

# ---- Hunk 4: remove lastRenderedPageBreak from "Trophy 2" paragraph ----
$idx = Find-ParaIndex $d "Trophy 2"
$inner = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Trophy 2</w:t></w:r>' +
    '</w:p>'
Replace-ParasWithXml $d $idx $idx $inner

Write-Host "Done"
